$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Tổng công tại LONG XUYÊN
$ws.Range("B12").Value = 22

# Lương cơ bản tại LONG XUYÊN
$ws.Range("B13").Value = 3142857.142857143

# Tổng lương tại LONG XUYÊN
$ws.Range("B32").Value = 792857.1428571432

# Tổng lương tại HỆ THỐNG
$ws.Range("B34").Value = 792857.1428571432
